# "adding averages and more checks"
#
# Refresh the Training Dashboard: the LAST UPDATE column (I) moves from
# 08-Sep-2025 to 16-Sep-2025 for every training row, and the PERIOD TO
# EXPIRE column (H) is recalculated accordingly (each value drops by the
# same 8 days). Also simplify the title/header font so both use a bold,
# white font instead of the old mismatched bold/size-14 vs plain-bold
# fonts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# New PERIOD TO EXPIRE values for rows 3-31 (each is 8 less than before,
# matching the LAST UPDATE date moving 8 days forward to 16-Sep-2025).
$newPeriod = @{
    3  = 672
    4  = 672
    5  = 672
    6  = 673
    7  = 673
    8  = 675
    9  = 675
    10 = 675
    11 = 687
    12 = 687
    13 = 674
    14 = 686
    15 = 691
    16 = 688
    17 = 674
    18 = 687
    19 = 686
    20 = 689
    21 = 674
    22 = 691
    23 = 686
    24 = 687
    25 = 689
    26 = 688
    27 = 686
    28 = 689
    29 = 308
    30 = 323
    31 = 324
}

foreach ($row in 3..31) {
    $ws.Cells.Item($row, 8).Value = $newPeriod[$row]
    # Leading apostrophe forces plain text entry so "16-Sep-2025" is
    # stored as a literal string (matching column I's existing text
    # values) instead of being auto-converted into a date serial.
    $ws.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# Restyle the title (row 1) and header row (row 2): both now use a
# bold, white font (the old distinct "bold size 14" title font and the
# plain "bold" header font collapse into one bold+white font, and the
# title loses its larger 14pt size).
$titleRange = $ws.Range("A1")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 11
$titleRange.Font.Color = 16777215

$headerRange = $ws.Range("A2:K2")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
